$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" value -------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) --------------------
#     The "Mapping: RIM Mapping" / "Mapping: Spécification métier vers
#     l'extension ROR FinancialHelpType" columns were reordered (AK and AL
#     swapped places), including their header, every data row, and the
#     column widths that were sized to fit their (now swapped) content.
$wsElem = $wb.Worksheets.Item("Elements")

# Only rows 1, 3, 5 and 6 actually hold differing AK/AL content; rows 2 and 4
# are blank in both columns, so leave them untouched.
$rowsToSwap = @(1, 3, 5, 6)
foreach ($r in $rowsToSwap) {
    $akCell = $wsElem.Cells.Item($r, 37)
    $alCell = $wsElem.Cells.Item($r, 38)
    $akValue = $akCell.Value2
    $alValue = $alCell.Value2
    $akCell.Value2 = $alValue
    $alCell.Value2 = $akValue
}

# Column widths follow the content that now lives in each column.
$wsElem.Columns.Item(37).ColumnWidth = 73.16666666666667
$wsElem.Columns.Item(38).ColumnWidth = 24.166666666666668
